$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as literal text (avoids numeric auto-coercion),
# then restore the default "Normal" style so no stray formatting is left behind.
function Set-TextValue([string]$cellRef, [string]$text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "60.859.27"
Set-TextValue "E2" "  +0.43%  "
Set-TextValue "D3" "2.593.67"
Set-TextValue "E3" "  +0.23%  "
Set-TextValue "E4" "  +0.02%  "
Set-TextValue "D5" "522.96"
Set-TextValue "E5" "  +2.97%  "
Set-TextValue "D6" "154.23"
Set-TextValue "E6" "  +0.40%  "
Set-TextValue "E7" "  +0.06%  "
Set-TextValue "E8" "  +2.60%  "
Set-TextValue "E9" "  +1.84%  "
Set-TextValue "E10" "  +1.14%  "
Set-TextValue "E11" "  -0.50%  "
Set-TextValue "D13" "3.048.48"
Set-TextValue "E13" "  +0.27%  "
Set-TextValue "D14" "60.885.25"
Set-TextValue "E14" "  +0.58%  "
Set-TextValue "E15" "  +0.27%  "
Set-TextValue "E16" "  +0.03%  "
Set-TextValue "D17" "2.599.91"
Set-TextValue "E17" "  +0.24%  "
Set-TextValue "D18" "4.75"
Set-TextValue "E18" "  -0.63%  "
Set-TextValue "D19" "353.89"
Set-TextValue "E19" "  +2.32%  "
Set-TextValue "D20" "10.56"
Set-TextValue "E20" "  +1.11%  "
Set-TextValue "E21" "  +1.33%  "
Set-TextValue "E22" "  +0.29%  "
Set-TextValue "D23" "60.79"
Set-TextValue "E23" "  +1.37%  "
Set-TextValue "E24" "  +1.42%  "
Set-TextValue "E25" "  -0.14%  "
Set-TextValue "D26" "2.708.62"
Set-TextValue "E26" "  +0.22%  "
Set-TextValue "D27" "1.00"
Set-TextValue "E27" "  +0.13%  "
Set-TextValue "E28" "  +0.13%  "
Set-TextValue "E29" "  +0.00%  "
Set-TextValue "E30" "  +0.01%  "
Set-TextValue "D31" "6.34"
Set-TextValue "E31" "  +10.81%  "
Set-TextValue "D32" "19.36"
Set-TextValue "E32" "  +0.06%  "
Set-TextValue "E33" "  +2.48%  "
Set-TextValue "D34" "148.13"
Set-TextValue "E34" "  -3.68%  "
Set-TextValue "E35" "  +4.34%  "
Set-TextValue "E36" "  +9.55%  "
Set-TextValue "E37" "  +0.91%  "
Set-TextValue "E38" "  +1.01%  "
Set-TextValue "E39" "  +1.86%  "
Set-TextValue "E40" "  +1.38%  "
Set-TextValue "E41" "  +1.36%  "
Set-TextValue "D42" "288.10"
Set-TextValue "E42" "  -2.26%  "
Set-TextValue "D43" "0.102"
Set-TextValue "E43" "  +2.04%  "
Set-TextValue "E44" "  -0.49%  "
Set-TextValue "D45" "0.0561"
Set-TextValue "E45" "  +0.64%  "
Set-TextValue "D46" "0.997"
Set-TextValue "E46" "  +0.06%  "
Set-TextValue "D47" "19.63"
Set-TextValue "E47" "  -0.98%  "
Set-TextValue "E48" "  +0.38%  "
Set-TextValue "E49" "  +1.99%  "
Set-TextValue "D50" "10.33"
Set-TextValue "E50" "  +0.34%  "
Set-TextValue "D51" "19.13"
Set-TextValue "E51" "  +8.20%  "
